$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 31

# Writing "2012-06-03" straight into BF2:BF31 would make Excel's
# autodetection re-interpret the text as a date serial, which would
# also pull a new number format into the stylesheet. Build the exact
# literal text with a helper formula, then paste back just the
# resulting values so the cells stay plain, unstyled text cells
# (same as the original "6-3-2011-12" inline strings).
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ws.Range("ZZ$r").Formula = '="2012-06-03"'
}

$helper = $ws.Range("ZZ$firstDataRow`:ZZ$lastDataRow")
$helper.Copy() | Out-Null

$target = $ws.Range("BF$firstDataRow`:BF$lastDataRow")
$target.PasteSpecial(-4163) | Out-Null  # xlPasteValues

$helper.ClearContents() | Out-Null
$excel.CutCopyMode = 0
